$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create sheets 2-4 as copies of sheet1 so they inherit the same sheetPr / ---
# --- sheetFormatPr / pageMargins boilerplate, then wipe + refill their content. ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws1.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Name = "Sales vs PO"
$ws2.Name = "Weekly Growth"
$ws3.Name = "Volume Insights"
$ws4.Name = "Prediction Info"

# --- Sheet1 "Sales vs PO": insert a new column C (shifts old PO_Requested_Qty ---
# --- column from C to D), rewrite column A (+6 days) and new column C/D.      ---
$ws1.Columns.Item(3).Insert()
$ws1.Cells.Item(1,3).Value = "Order Week"
$ws1.Cells.Item(2,1).Value = 45445
$ws1.Cells.Item(2,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(2,3).Value = 45439
$ws1.Cells.Item(2,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(2,4).Value = 0
$ws1.Cells.Item(3,1).Value = 45452
$ws1.Cells.Item(3,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(3,3).Value = 45446
$ws1.Cells.Item(3,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(3,4).Value = 0
$ws1.Cells.Item(4,1).Value = 45459
$ws1.Cells.Item(4,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(4,3).Value = 45453
$ws1.Cells.Item(4,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(4,4).Value = 0
$ws1.Cells.Item(5,1).Value = 45466
$ws1.Cells.Item(5,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(5,3).Value = 45460
$ws1.Cells.Item(5,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(5,4).Value = 0
$ws1.Cells.Item(6,1).Value = 45473
$ws1.Cells.Item(6,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(6,3).Value = 45467
$ws1.Cells.Item(6,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(6,4).Value = 0
$ws1.Cells.Item(7,1).Value = 45480
$ws1.Cells.Item(7,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(7,3).Value = 45474
$ws1.Cells.Item(7,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(7,4).Value = 0
$ws1.Cells.Item(8,1).Value = 45487
$ws1.Cells.Item(8,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(8,3).Value = 45481
$ws1.Cells.Item(8,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(8,4).Value = 0
$ws1.Cells.Item(9,1).Value = 45494
$ws1.Cells.Item(9,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(9,3).Value = 45488
$ws1.Cells.Item(9,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(9,4).Value = 0
$ws1.Cells.Item(10,1).Value = 45501
$ws1.Cells.Item(10,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(10,3).Value = 45495
$ws1.Cells.Item(10,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(10,4).Value = 0
$ws1.Cells.Item(11,1).Value = 45508
$ws1.Cells.Item(11,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(11,3).Value = 45502
$ws1.Cells.Item(11,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(11,4).Value = 0
$ws1.Cells.Item(12,1).Value = 45515
$ws1.Cells.Item(12,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(12,3).Value = 45509
$ws1.Cells.Item(12,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(12,4).Value = 0
$ws1.Cells.Item(13,1).Value = 45522
$ws1.Cells.Item(13,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(13,3).Value = 45516
$ws1.Cells.Item(13,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(13,4).Value = 0
$ws1.Cells.Item(14,1).Value = 45529
$ws1.Cells.Item(14,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(14,3).Value = 45523
$ws1.Cells.Item(14,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(14,4).Value = 0
$ws1.Cells.Item(15,1).Value = 45536
$ws1.Cells.Item(15,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(15,3).Value = 45530
$ws1.Cells.Item(15,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(15,4).Value = 0
$ws1.Cells.Item(16,1).Value = 45543
$ws1.Cells.Item(16,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(16,3).Value = 45537
$ws1.Cells.Item(16,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(16,4).Value = 0
$ws1.Cells.Item(17,1).Value = 45550
$ws1.Cells.Item(17,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(17,3).Value = 45544
$ws1.Cells.Item(17,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(17,4).Value = 0
$ws1.Cells.Item(18,1).Value = 45557
$ws1.Cells.Item(18,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(18,3).Value = 45551
$ws1.Cells.Item(18,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(18,4).Value = 0
$ws1.Cells.Item(19,1).Value = 45564
$ws1.Cells.Item(19,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(19,3).Value = 45558
$ws1.Cells.Item(19,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(19,4).Value = 0
$ws1.Cells.Item(20,1).Value = 45571
$ws1.Cells.Item(20,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(20,3).Value = 45565
$ws1.Cells.Item(20,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(20,4).Value = 0
$ws1.Cells.Item(21,1).Value = 45578
$ws1.Cells.Item(21,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(21,3).Value = 45572
$ws1.Cells.Item(21,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(21,4).Value = 0
$ws1.Cells.Item(22,1).Value = 45585
$ws1.Cells.Item(22,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(22,3).Value = 45579
$ws1.Cells.Item(22,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(22,4).Value = 0
$ws1.Cells.Item(23,1).Value = 45592
$ws1.Cells.Item(23,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(23,3).Value = 45586
$ws1.Cells.Item(23,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(23,4).Value = 0
$ws1.Cells.Item(24,1).Value = 45599
$ws1.Cells.Item(24,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(24,3).Value = 45593
$ws1.Cells.Item(24,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(24,4).Value = 0
$ws1.Cells.Item(25,1).Value = 45606
$ws1.Cells.Item(25,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(25,3).Value = 45600
$ws1.Cells.Item(25,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(25,4).Value = 0
$ws1.Cells.Item(26,1).Value = 45613
$ws1.Cells.Item(26,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(26,3).Value = 45607
$ws1.Cells.Item(26,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(26,4).Value = 0
$ws1.Cells.Item(27,1).Value = 45620
$ws1.Cells.Item(27,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(27,3).Value = 45614
$ws1.Cells.Item(27,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(27,4).Value = 0
$ws1.Cells.Item(28,1).Value = 45627
$ws1.Cells.Item(28,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(28,3).Value = 45621
$ws1.Cells.Item(28,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(28,4).Value = 0
$ws1.Cells.Item(29,1).Value = 45634
$ws1.Cells.Item(29,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(29,3).Value = 45628
$ws1.Cells.Item(29,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(29,4).Value = 0
$ws1.Cells.Item(30,1).Value = 45641
$ws1.Cells.Item(30,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(30,3).Value = 45635
$ws1.Cells.Item(30,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(30,4).Value = 0
$ws1.Cells.Item(31,1).Value = 45648
$ws1.Cells.Item(31,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(31,3).Value = 45642
$ws1.Cells.Item(31,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(31,4).Value = 0
$ws1.Cells.Item(32,1).Value = 45655
$ws1.Cells.Item(32,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(32,3).Value = 45649
$ws1.Cells.Item(32,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(32,4).Value = 0

# --- Sheet2: Weekly Growth ---
$ws2.Cells.Clear()
$ws2.Cells.Item(1,1).Value = "ds"
$ws2.Cells.Item(1,2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1,3).Value = "Growth%"
$ws1.Range("A1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws2.Cells.Item(2,1).Value = 45446
$ws2.Cells.Item(2,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(2,2).Value = 180
$ws2.Cells.Item(2,3).Value = 0
$ws2.Cells.Item(3,1).Value = 45460
$ws2.Cells.Item(3,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(3,2).Value = 400
$ws2.Cells.Item(3,3).Value = 122.2222222222222
$ws2.Cells.Item(4,1).Value = 45467
$ws2.Cells.Item(4,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(4,2).Value = 100
$ws2.Cells.Item(4,3).Value = -75
$ws2.Cells.Item(5,1).Value = 45474
$ws2.Cells.Item(5,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(5,2).Value = 80
$ws2.Cells.Item(5,3).Value = -20
$ws2.Cells.Item(6,1).Value = 45481
$ws2.Cells.Item(6,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(6,2).Value = 40
$ws2.Cells.Item(6,3).Value = -50
$ws2.Cells.Item(7,1).Value = 45488
$ws2.Cells.Item(7,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(7,2).Value = 260
$ws2.Cells.Item(7,3).Value = 550
$ws2.Cells.Item(8,1).Value = 45516
$ws2.Cells.Item(8,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(8,2).Value = 100
$ws2.Cells.Item(8,3).Value = -61.53846153846154
$ws2.Cells.Item(9,1).Value = 45523
$ws2.Cells.Item(9,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(9,2).Value = 160
$ws2.Cells.Item(9,3).Value = 60.00000000000001
$ws2.Cells.Item(10,1).Value = 45530
$ws2.Cells.Item(10,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(10,2).Value = 80
$ws2.Cells.Item(10,3).Value = -50
$ws2.Cells.Item(11,1).Value = 45537
$ws2.Cells.Item(11,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(11,2).Value = 120
$ws2.Cells.Item(11,3).Value = 50
$ws2.Cells.Item(12,1).Value = 45544
$ws2.Cells.Item(12,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(12,2).Value = 40
$ws2.Cells.Item(12,3).Value = -66.66666666666667
$ws2.Cells.Item(13,1).Value = 45551
$ws2.Cells.Item(13,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(13,2).Value = 280
$ws2.Cells.Item(13,3).Value = 600
$ws2.Cells.Item(14,1).Value = 45572
$ws2.Cells.Item(14,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(14,2).Value = 60
$ws2.Cells.Item(14,3).Value = -78.57142857142857
$ws2.Cells.Item(15,1).Value = 45579
$ws2.Cells.Item(15,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(15,2).Value = 20
$ws2.Cells.Item(15,3).Value = -66.66666666666667
$ws2.Cells.Item(16,1).Value = 45586
$ws2.Cells.Item(16,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(16,2).Value = 20
$ws2.Cells.Item(16,3).Value = 0
$ws2.Cells.Item(17,1).Value = 45593
$ws2.Cells.Item(17,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(17,2).Value = 20
$ws2.Cells.Item(17,3).Value = 0
$ws2.Cells.Item(18,1).Value = 45600
$ws2.Cells.Item(18,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(18,2).Value = 20
$ws2.Cells.Item(18,3).Value = 0
$ws2.Cells.Item(19,1).Value = 45628
$ws2.Cells.Item(19,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(19,2).Value = 20
$ws2.Cells.Item(19,3).Value = 0
$ws2.Cells.Item(20,1).Value = 45635
$ws2.Cells.Item(20,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(20,2).Value = 20
$ws2.Cells.Item(20,3).Value = 0

# --- Sheet3: Volume Insights ---
$ws3.Cells.Clear()
$ws3.Cells.Item(1,1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1,2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1,3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1,4).Value = "Min_PO_Quantity"
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws3.Cells.Item(2,1).Value = 2020
$ws3.Cells.Item(2,2).Value = 106.3157894736842
$ws3.Cells.Item(2,3).Value = 400
$ws3.Cells.Item(2,4).Value = 20

# --- Sheet4: Prediction Info ---
$ws4.Cells.Clear()
$ws4.Cells.Item(1,1).Value = "Predicted_Next_Week_PO_Quantity"
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$ws4.Cells.Item(2,1).Value = 0

$ws1.Select()
$ws1.Range("A1").Select()
